# Izraden js za polja 'datum od' i 'datum do' da se ravnaju po polju 'na dan'.
# Ovdje: uklanjamo stupce 'datum_od' i 'datum_do' (A:B) iz lista 'Zaglavlje',
# jer se ti datumi sada automatski racunaju iz polja 'na dan'.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zaglavlje")

# Brisanje stupaca A (datum_od) i B (datum_do) - ostali stupci se pomicu ulijevo.
$ws.Range("A1:B1").EntireColumn.Delete()

# Uskladi odabir prikaza sa novim rasporedom stupaca.
$ws.Range("E9").Select()
